$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20:101 down to 21:102
$ws.Rows("20:20").Insert()

# Populate the new row 20 with the new data record
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 45063
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107011
$ws.Range("J20").Value = "Tuna"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = "$/caja 16 kilos"
$ws.Range("R20").Value = "Provincia de Los Andes"
$ws.Range("S20").Value = 1250
$ws.Range("T20").Value = 16
